$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 12.59759837740221
$ws.Range("D2").Value = 12.07807400996011
$ws.Range("E2").Value = 12.38685093838693
$ws.Range("F2").Value = 59.54525487882027
$ws.Range("G2").Value = 3.791290863414195
$ws.Range("K2").Value = 22.2946968890948
$ws.Range("L2").Value = 8.620385607774349
$ws.Range("M2").Value = 22.56752991014817

$ws.Range("C3").Value = 12.56548497165073
$ws.Range("D3").Value = 11.99063217526515
$ws.Range("E3").Value = 12.3918448796539
$ws.Range("F3").Value = 58.64943377490314
$ws.Range("G3").Value = 3.796486930794841
$ws.Range("K3").Value = 22.17715874727757
$ws.Range("L3").Value = 8.633470237417729
$ws.Range("M3").Value = 22.54860712892781

$ws.Range("C4").Value = 12.54849518781128
$ws.Range("D4").Value = 11.93699011476378
$ws.Range("E4").Value = 12.39671666261354
$ws.Range("F4").Value = 58.1010366047169
$ws.Range("G4").Value = 3.799835133562097
$ws.Range("K4").Value = 22.11375484052347
$ws.Range("L4").Value = 8.642112729878654
$ws.Range("M4").Value = 22.54381306196882

$ws.Range("C5").Value = 12.54225972679078
$ws.Range("D5").Value = 11.91514738260152
$ws.Range("E5").Value = 12.39915520495384
$ws.Range("F5").Value = 57.87816530659294
$ws.Range("G5").Value = 3.801239421090784
$ws.Range("K5").Value = 22.09014256364203
$ws.Range("L5").Value = 8.645787852764725
$ws.Range("M5").Value = 22.5435749621522

$ws.Range("C6").Value = 12.54126593447682
$ws.Range("D6").Value = 11.91152156771346
$ws.Range("E6").Value = 12.39958747427976
$ws.Range("F6").Value = 57.84119983154423
$ws.Range("G6").Value = 3.801475015394334
$ws.Range("K6").Value = 22.08635670965693
$ws.Range("L6").Value = 8.646407366077016
$ws.Range("M6").Value = 22.54363899459797

$ws.Range("C7").Value = 12.54840830622855
$ws.Range("D7").Value = 11.93669546161947
$ws.Range("E7").Value = 12.39674771550861
$ws.Range("F7").Value = 58.09802818164776
$ws.Range("G7").Value = 3.799853910600621
$ws.Range("K7").Value = 22.11342736164215
$ws.Range("L7").Value = 8.642161673106816
$ws.Range("M7").Value = 22.54380290654377

$ws.Range("C8").Value = 12.58596089707112
$ws.Range("D8").Value = 12.04791718998475
$ws.Range("E8").Value = 12.38819769595746
$ws.Range("F8").Value = 59.23615101768195
$ws.Range("G8").Value = 3.793049830700858
$ws.Range("K8").Value = 22.25236458316238
$ws.Range("L8").Value = 8.624771036817629
$ws.Range("M8").Value = 22.55958895413655

$ws.Range("C9").Value = 12.6811304610611
$ws.Range("D9").Value = 12.26617882926515
$ws.Range("E9").Value = 12.3857880830086
$ws.Range("F9").Value = 61.47256210335246
$ws.Range("G9").Value = 3.780950473889391
$ws.Range("K9").Value = 22.59323674705218
$ws.Range("L9").Value = 8.595485869538527
$ws.Range("M9").Value = 22.64466153940324

$ws.Range("C10").Value = 12.76398281787807
$ws.Range("D10").Value = 12.42632645371046
$ws.Range("E10").Value = 12.39280945162262
$ws.Range("F10").Value = 63.10719892436272
$ws.Range("G10").Value = 3.772807001868363
$ws.Range("K10").Value = 22.88360863703198
$ws.Range("L10").Value = 8.576893393334066
$ws.Range("M10").Value = 22.74002287599276

$ws.Range("C11").Value = 12.80442823018323
$ws.Range("D11").Value = 12.49906757210049
$ws.Range("E11").Value = 12.39791899303138
$ws.Range("F11").Value = 63.84667549704835
$ws.Range("G11").Value = 3.769261689641248
$ws.Range("K11").Value = 23.02392692633866
$ws.Range("L11").Value = 8.569067125122604
$ws.Range("M11").Value = 22.79048197367107

$ws.Range("C12").Value = 12.82013445380957
$ws.Range("D12").Value = 12.52659052574553
$ws.Range("E12").Value = 12.40012954378843
$ws.Range("F12").Value = 64.12591340395619
$ws.Range("G12").Value = 3.767941860313335
$ws.Range("K12").Value = 23.07820245830274
$ws.Range("L12").Value = 8.566194122043713
$ws.Range("M12").Value = 22.81060009139166

$ws.Range("C13").Value = 12.81673457198417
$ws.Range("D13").Value = 12.52066407351463
$ws.Range("E13").Value = 12.39964119882293
$ws.Range("F13").Value = 64.06581246924543
$ws.Range("G13").Value = 3.768225102435389
$ws.Range("K13").Value = 23.06646325838554
$ws.Range("L13").Value = 8.5668088464804
$ws.Range("M13").Value = 22.80622249456136

$ws.Range("C14").Value = 12.80571260130553
$ws.Range("D14").Value = 12.50133239176449
$ws.Range("E14").Value = 12.39809533045392
$ws.Range("F14").Value = 63.8696653960135
$ws.Range("G14").Value = 3.7691526525584
$ws.Range("K14").Value = 23.02836961392979
$ws.Range("L14").Value = 8.568828945991596
$ws.Range("M14").Value = 22.79211690273297

$ws.Range("C15").Value = 12.79901198842334
$ws.Range("D15").Value = 12.48948804680185
$ws.Range("E15").Value = 12.39718434779406
$ws.Range("F15").Value = 63.749411711731
$ws.Range("G15").Value = 3.769723754850571
$ws.Range("K15").Value = 23.00518330552013
$ws.Range("L15").Value = 8.570078113320871
$ws.Range("M15").Value = 22.78360814940987

$ws.Range("C16").Value = 12.76139457885994
$ws.Range("D16").Value = 12.42156985979878
$ws.Range("E16").Value = 12.39251409826884
$ws.Range("F16").Value = 63.05877275796436
$ws.Range("G16").Value = 3.773041883654264
$ws.Range("K16").Value = 22.874600437128
$ws.Range("L16").Value = 8.577417551124626
$ws.Range("M16").Value = 22.73686708391963

$ws.Range("C17").Value = 12.73901931061277
$ws.Range("D17").Value = 12.37987121710253
$ws.Range("E17").Value = 12.39013985767228
$ws.Range("F17").Value = 62.63389223562348
$ws.Range("G17").Value = 3.775118087677429
$ws.Range("K17").Value = 22.79656941238293
$ws.Range("L17").Value = 8.582081682224887
$ws.Range("M17").Value = 22.71000106516975

$ws.Range("C18").Value = 12.72640939088515
$ws.Range("D18").Value = 12.35587723180245
$ws.Range("E18").Value = 12.38895454553175
$ws.Range("F18").Value = 62.3891348123399
$ws.Range("G18").Value = 3.776327260714488
$ws.Range("K18").Value = 22.75246511269575
$ws.Range("L18").Value = 8.584823821389035
$ws.Range("M18").Value = 22.69521527185261

$ws.Range("C19").Value = 12.72218466713701
$ws.Range("D19").Value = 12.34775175686172
$ws.Range("E19").Value = 12.38858417479306
$ws.Range("F19").Value = 62.30620529234125
$ws.Range("G19").Value = 3.776739247131822
$ws.Range("K19").Value = 22.73766688424997
$ws.Range("L19").Value = 8.585762479476783
$ws.Range("M19").Value = 22.69032379187032

$ws.Range("C20").Value = 12.74137435394368
$ws.Range("D20").Value = 12.3843111940796
$ws.Range("E20").Value = 12.39037393841677
$ws.Range("F20").Value = 62.6791619620344
$ws.Range("G20").Value = 3.774895521612773
$ws.Range("K20").Value = 22.80479582598167
$ws.Range("L20").Value = 8.581579025530274
$ws.Range("M20").Value = 22.71279203461697

$ws.Range("C21").Value = 12.8089394744473
$ws.Range("D21").Value = 12.50701123988815
$ws.Range("E21").Value = 12.39854190636838
$ws.Range("F21").Value = 63.92730137461005
$ws.Range("G21").Value = 3.768879593941048
$ws.Range("K21").Value = 23.03952805516397
$ws.Range("L21").Value = 8.568233135410761
$ws.Range("M21").Value = 22.79623270360759

$ws.Range("C22").Value = 12.85536977815451
$ws.Range("D22").Value = 12.58706835751228
$ws.Range("E22").Value = 12.40548703372143
$ws.Range("F22").Value = 64.73836761101064
$ws.Range("G22").Value = 3.765080082800936
$ws.Range("K22").Value = 23.19956253281338
$ws.Range("L22").Value = 8.560039033067708
$ws.Range("M22").Value = 22.85665060976613

$ws.Range("C23").Value = 12.83038317608532
$ws.Range("D23").Value = 12.54435487785539
$ws.Range("E23").Value = 12.4016332118703
$ws.Range("F23").Value = 64.30597502495847
$ws.Range("G23").Value = 3.767095915335587
$ws.Range("K23").Value = 23.11355790542418
$ws.Range("L23").Value = 8.564364107121323
$ws.Range("M23").Value = 22.82386885786406

$ws.Range("C24").Value = 12.74030884751613
$ws.Range("D24").Value = 12.3823039450264
$ws.Range("E24").Value = 12.3902675508981
$ws.Range("F24").Value = 62.65869702393449
$ws.Range("G24").Value = 3.7749960953263
$ws.Range("K24").Value = 22.80107430397972
$ws.Range("L24").Value = 8.581806087640844
$ws.Range("M24").Value = 22.71152818152317

$ws.Range("C25").Value = 12.65309715373743
$ws.Range("D25").Value = 12.20714706364791
$ws.Range("E25").Value = 12.38489842350936
$ws.Range("F25").Value = 60.86823420918584
$ws.Range("G25").Value = 3.784091802996222
$ws.Range("K25").Value = 22.49386643115766
$ws.Range("L25").Value = 8.602893935878173
$ws.Range("M25").Value = 22.61586150137859
